$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "ueiDUNS" columns that are no longer part of the export
# (hqParent.ueiDUNS and ultimateParentEntity.ueiDUNS). Deleting BB first
# keeps the AQ reference valid for the second delete.
$ws.Columns("BB").Delete()
$ws.Columns("AQ").Delete()

$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
